# aggiornamento fino a 02/05
# Append rolling-7gg municipality counts for 2021-04-27 .. 2021-05-02
# as new rows 239-244 (columns A:AX), extending the existing series.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 6,50
$rowVals = @(44313,29,7,10,68,57,7,21,6,11,1,9,32,43,1,2,0,55,6,9,28,249,9,13,13,45,4,11,2,6,19,6,94,10,10,2,14,17,61,7,30,1045,10,0,1,2,3,2,3,0)
for ($c = 0; $c -lt 50; $c++) { $data[0,$c] = $rowVals[$c] }
$rowVals = @(44314,23,7,8,68,53,8,21,6,12,1,9,27,43,1,2,0,49,6,9,29,258,10,15,13,46,2,7,2,5,19,6,82,10,10,1,13,21,59,7,27,1015,10,0,1,2,3,1,3,0)
for ($c = 0; $c -lt 50; $c++) { $data[1,$c] = $rowVals[$c] }
$rowVals = @(44315,24,9,6,71,54,7,21,6,13,2,14,28,42,0,2,0,43,5,11,28,237,15,19,14,42,1,7,1,5,13,9,78,14,8,0,14,22,58,4,25,989,10,0,1,2,0,1,3,0)
for ($c = 0; $c -lt 50; $c++) { $data[2,$c] = $rowVals[$c] }
$rowVals = @(44316,31,8,6,92,44,8,21,9,13,1,13,33,44,0,3,1,31,5,12,27,241,19,16,8,40,1,8,0,5,11,10,67,21,7,1,18,26,54,1,23,995,10,0,1,3,0,0,2,0)
for ($c = 0; $c -lt 50; $c++) { $data[3,$c] = $rowVals[$c] }
$rowVals = @(44317,34,6,5,94,43,12,18,11,14,2,14,28,44,0,6,1,37,6,11,29,255,21,20,8,34,1,4,0,3,10,13,71,26,8,1,20,30,53,1,28,1040,9,1,1,5,0,1,1,0)
for ($c = 0; $c -lt 50; $c++) { $data[4,$c] = $rowVals[$c] }
$rowVals = @(44318,34,4,1,94,39,10,20,13,12,3,14,33,45,0,5,2,41,7,9,26,247,20,19,9,28,2,4,0,4,12,11,83,28,12,1,19,25,38,1,29,1017,3,1,1,5,1,1,1,0)
for ($c = 0; $c -lt 50; $c++) { $data[5,$c] = $rowVals[$c] }

$ws.Range("A239:AX244").Value = $data

# Column A (date) reuses the bold / centered / bordered date style
# already used by the rest of column A (copied from the prior last row).
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)
$excel.CutCopyMode = $false
